$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple "Taxonsorteringsordning" (column B) value bumps ---
$ws.Range("B2").Value  = 79243
$ws.Range("B3").Value  = 79244
$ws.Range("B4").Value  = 91804
$ws.Range("B5").Value  = 79243
$ws.Range("B6").Value  = 79243
$ws.Range("B7").Value  = 91804
$ws.Range("B8").Value  = 91828
$ws.Range("B9").Value  = 79243
$ws.Range("B12").Value = 91828
$ws.Range("B13").Value = 79243
$ws.Range("B14").Value = 91804
$ws.Range("B15").Value = 57884
$ws.Range("B16").Value = 79243
$ws.Range("B17").Value = 79243
$ws.Range("B18").Value = 79243
$ws.Range("B19").Value = 79243
$ws.Range("B20").Value = 57884
$ws.Range("B21").Value = 57884

# --- Rows 10 and 11 swap their record content (with updated B values) ---

# Row 10 becomes the "Tretåig hackspett" record (previously in row 11)
$ws.Range("A10").Value = 130861151
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = "Tretåig hackspett"
$ws.Range("G10").Value = "Picoides tridactylus"
$ws.Range("H10").Value = "(Linnaeus, 1758)"
$ws.Range("M10").Value = "äldre spår"
$ws.Range("Q10").Value = 442749
$ws.Range("R10").Value = 7039568
$ws.Range("AC10").Value = "Ringhack, äldre, enstaka på en gran. Mycket högt livsmiljövärde för tretåig hackspett kring fyndplatsen."
$ws.Range("AH10").Value = "Granskog"
$ws.Range("AJ10").Value = "gran"
$ws.Range("AK10").Value = "Picea abies"
$ws.Range("AM10").Value = "Trädstam på levande träd"
$ws.Range("AO10").Value = "Stem on living tree # Picea abies"

# Row 11 becomes the "Garnlav" record (previously in row 10)
$ws.Range("A11").Value = 130861155
$ws.Range("B11").Value = 79243
$ws.Range("E11").Value = 6425
$ws.Range("F11").Value = "Garnlav"
$ws.Range("G11").Value = "Alectoria sarmentosa"
$ws.Range("H11").Value = "(Ach.) Ach."
$ws.Range("M11").ClearContents()
$ws.Range("Q11").Value = 442870
$ws.Range("R11").Value = 7039632
$ws.Range("AC11").Value = "På död undertryck gran."
$ws.Range("AH11").ClearContents()
$ws.Range("AJ11").ClearContents()
$ws.Range("AK11").ClearContents()
$ws.Range("AM11").ClearContents()
$ws.Range("AO11").ClearContents()
